# Update the "K" column (column G) values for rows 2-8 on the active sheet.
# These values were regenerated upstream (K computed from strikeouts instead
# of the old "Strike#" metric) - write the new literal results here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
